$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = @(-0.37990752884762458, -0.28076836278126827, -0.17526878804605772, -0.16255050431589524, -0.15590410052867121, -0.053693907042594358, -0.0336512189971927, -0.013590046044767767, -0.0075767881138117943, -0.0015730248552827675, 0.0029252738571514669, 0.0089244961061694106, 0.014936307395704951, 0.026947405111591216, -0.021916253928893958, -0.01589265180273447, -0.0098680625889997486, -0.088291782639103644, -0.07921591797002403, -0.01801312625929441, -0.0090041179148014905, -0.093948180456331798, -0.084635256003168458, -0.042126821611276455, -0.094917072273020153, -0.088791059760069402, -0.083658898681844729, -0.077047264706932239, -0.064643955675490972, -0.044466820664855078, -0.029279609117857319, -0.0060004942579583442)
$colB = @(0.37906014479179362, 0.27821696786047845, 0.17455050414512741, 0.16190409989125243, 0.15461097865715612, 0.053651218226137942, 0.033590045270232238, 0.013576787448246641, 0.0075730241885985095, 0.0015747254880338346, -0.0029244967727990456, -0.0089363080601589928, -0.014947405822335114, -0.026987707487334589, 0.021892651139609143, 0.015868061922714727, 0.009861174282823626, 0.088215917347582717, 0.078622992151370585, 0.018004117273063613, 0.0089999993574343407, 0.093635255355657421, 0.084126820690291204, 0.041999999073254024, 0.094791059099208042, 0.088628573342880657, 0.083047264027666934, 0.076643954938231396, 0.064466819858715496, 0.044279608343875676, 0.029235415611902127, 0.0059999992958399417)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $colA[$i]
    $ws.Cells.Item($i + 1, 2).Value = $colB[$i]
}

# Column B width changes from 14.7109375 to 15.42578125 (same as column A).
# The COM ColumnWidth property snaps to a 1/6-character pixel grid, so we use the
# input value whose rounded result lands closest to the target stored width.
$ws.Columns.Item(2).ColumnWidth = 14.66
